# Generate Report for handoff
#
# The file "5fe91086-255c-4682-9774-a18afabf88ac.md" has just been handed
# off (for both zh-cn and de-de targets), so its status flips from
# "Handed back" to "Not yet handed off", and the corresponding "Latest
# Handoff Datetime" cells are stamped with the new handoff time.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: summary status for the file, one column per locale ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Not yet handed off"
$overview.Range("C3").Value = "Not yet handed off"

# --- zh-cn sheet: detailed row for the file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Not yet handed off"
$zhcn.Range("D3").Value = "2016-01-08 14:17:44"

# --- de-de sheet: detailed row for the file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Not yet handed off"
$dede.Range("D3").Value = "2016-01-08 14:17:59"
